# ER atualizado e tabelas mapeadas na pasta 'tabelas novas'
#
# The "hora_inicio/hour_begin" and "hora_fim/hour_end" rows were merged
# into the "data_inicio" / "data_fim" rows: the new-name column now reads
# "date_hour_begin" / "date_hour_end" and the separate hour rows were
# cleared out (leaving blank spacer rows in their place).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5: data_inicio -> date_hour_begin (was date_begin)
$ws.Range("B5").Value = "date_hour_begin"

# Row 6 used to hold hora_inicio / hour_begin - now blank
$ws.Range("A6").Value = ""
$ws.Range("B6").Value = ""

# Row 7: data_fim -> date_hour_end (was date_end)
$ws.Range("B7").Value = "date_hour_end"

# Row 8 used to hold hora_fim / hour_end - now blank
$ws.Range("A8").Value = ""
$ws.Range("B8").Value = ""

$ws.Range("A7").Select()
